$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Lesson plan content updates -------------------------------------------------

# Topic cells for week03/week04 (row3) got swapped/updated, and week04 (row4)
# "Mi az az internet?" became "Mi az az internet 2?".
$ws.Range("B3").Value = "Mi az az internet?"
$ws.Range("C3").Value = "Mi az a wifi?"
$ws.Range("B4").Value = "Mi az az internet 2?"

# Fix stray trailing character typo in the wifi-frequency bullet point.
$ws.Range("J8").Value = "különböző frekvencia(fontos hogy ne legyen interferencia, pl emergency services)"

# The first two "Lesson.. - .." header labels on top of the J/K list box were
# renumbered (Lesson01 -> Lesson03, Lesson02 -> Lesson0102).
$ws.Range("J2").Value = "Lesson03 - Mi az a wifi?"
$ws.Range("K2").Value = "Lesson0102 - Mi az az internet?"

# --- Extend the bordered J:K list box down by one row -----------------------------

# K12 currently carries the box's thick bottom border (style with bottom edge);
# push that border down onto the new K13 cell, then give K12 the plain
# (no-bottom-border) interior style used by the rest of the box.
$ws.Range("K12").Copy()
$ws.Range("K13").PasteSpecial(-4122)

$ws.Range("K3").Copy()
$ws.Range("K12").PasteSpecial(-4122)

# Row 12 no longer needs the taller "thick bottom border" row height now that
# its border moved to row 13 - auto-fit it back down to the standard height.
$ws.Rows.Item(12).AutoFit()

# --- Selection / view state --------------------------------------------------------

$ws.Range("J2:K13").Select()

# --- Column widths -------------------------------------------------------------

# Column B and K hold longer text now (best-fit widths need to grow).
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(11).AutoFit()
